# The sheet had its data sitting in columns B:F with an unused/blank column A.
# The edit removes that blank leading column so the table now starts at column A
# (B:F -> A:E), and fills in the two previously-empty cells in the new D:E
# columns of row 2 with 0 (they were blank before the shift).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the whole table one column to the left by deleting the empty column A.
$ws.Range("A:A").Delete()

# Row 2 only had 3 values before (now in A2:C2); the shift exposes new cells
# D2:E2 that Excel leaves blank. Populate them with 0 to match the data set.
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

# Restore/update the active selection as left by the editor.
$ws.Range("G5").Select() | Out-Null
